$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.185.03'
$ws.Range('E2').Value = '  +2.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.596.24'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '531.73'
$ws.Range('E5').Value = '  +3.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.90'
$ws.Range('E6').Value = '  +1.96%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +1.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.609.65'
$ws.Range('E9').Value = '  +1.10%  '
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('E11').Value = '  +3.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.335'
$ws.Range('E12').Value = '  +3.24%  '
$ws.Range('E13').Value = '  +2.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.056.52'
$ws.Range('E14').Value = '  +1.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.133.95'
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.49'
$ws.Range('E16').Value = '  +2.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.609.81'
$ws.Range('E17').Value = '  +1.55%  '
$ws.Range('E18').Value = '  +2.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '346.74'
$ws.Range('E19').Value = '  +4.76%  '
$ws.Range('E20').Value = '  +1.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.13'
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.39'
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.46'
$ws.Range('E24').Value = '  +2.50%  '
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.406'
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.16'
$ws.Range('E28').Value = '  +3.99%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0734'
$ws.Range('E30').Value = '  +3.45%  '
$ws.Range('E31').Value = '  +4.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.82'
$ws.Range('E32').Value = '  -1.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.78'
$ws.Range('E33').Value = '  +1.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.53'
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.97'
$ws.Range('E35').Value = '  +2.21%  '
$ws.Range('E36').Value = '  +1.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '36.83'
$ws.Range('E37').Value = '  +2.06%  '
$ws.Range('E38').Value = '  +4.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.838'
$ws.Range('E39').Value = '  +2.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.840'
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.53'
$ws.Range('E41').Value = '  +1.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.997'
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '272.31'
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('E45').Value = '  +1.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0960'
$ws.Range('E46').Value = '  +2.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0520'
$ws.Range('E47').Value = '  +1.61%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.65'
$ws.Range('E48').Value = '  +5.69%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.947.33'
$ws.Range('E49').Value = '  -1.00%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0222'
$ws.Range('E50').Value = '  +2.48%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.28'
$ws.Range('E51').Value = '  +2.98%  '
